$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1586
$ws.Range("H112").Value = 2048.7446
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 2069.3696
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 6208.1088
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -8424.1088
$ws.Range("H138").Value = 822363.5600000001
$ws.Range("I138").Value = 1724.625
$ws.Range("J138").Value = 968254.9399999999
$ws.Range("K138").Value = 5173.875
$ws.Range("L138").Value = 2904764.82
$ws.Range("M138").Value = -33.875
$ws.Range("N138").Value = -2915044.82

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4544.067
$ws.Range("I32").Value = 4878
$ws.Range("K32").Value = 4878
$ws.Range("M32").Value = -4591
$ws.Range("H45").Value = 1624.6154
$ws.Range("I45").Value = 1624.6154
$ws.Range("K45").Value = 1624.6154
$ws.Range("M45").Value = -1247.6154

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4153
$ws.Range("I86").Value = 4575.4
$ws.Range("J86").Value = 3214.3333
$ws.Range("K86").Value = 4575.4
$ws.Range("L86").Value = 3214.3333
$ws.Range("M86").Value = -3452.4
$ws.Range("N86").Value = -5460.3333
$ws.Range("H89").Value = 4153
$ws.Range("I89").Value = 4575.4
$ws.Range("J89").Value = 3214.3333
$ws.Range("K89").Value = 22877
$ws.Range("L89").Value = 16071.6665
$ws.Range("M89").Value = -17261
$ws.Range("N89").Value = -27303.6665

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1932.0625
$ws.Range("I5").Value = 2064.8572
$ws.Range("K5").Value = 6194.571599999999
$ws.Range("M5").Value = -6082.571599999999
$ws.Range("H31").Value = 1201.5186
$ws.Range("I31").Value = 1057.64
$ws.Range("K31").Value = 1057.64
$ws.Range("M31").Value = -762.6400000000001
$ws.Range("H34").Value = 1201.5186
$ws.Range("I34").Value = 1057.64
$ws.Range("K34").Value = 1057.64
$ws.Range("M34").Value = -855.6400000000001
$ws.Range("H58").Value = 742.3333
$ws.Range("I58").Value = 725.36365
$ws.Range("J58").Value = 789
$ws.Range("K58").Value = 725.36365
$ws.Range("L58").Value = 789
$ws.Range("M58").Value = -522.36365
$ws.Range("N58").Value = -1195
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H115").Value = 59999
$ws.Range("J115").Value = 59999
$ws.Range("L115").Value = 59999
$ws.Range("N115").Value = -62349
$ws.Range("H136").Value = 742.3333
$ws.Range("I136").Value = 725.36365
$ws.Range("J136").Value = 789
$ws.Range("K136").Value = 2176.09095
$ws.Range("L136").Value = 2367
$ws.Range("M136").Value = 373.9090500000002
$ws.Range("N136").Value = -7467

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4181.727
$ws.Range("I64").Value = 799.5
$ws.Range("J64").Value = 4933.3335
$ws.Range("K64").Value = 2398.5
$ws.Range("L64").Value = 14800.0005
$ws.Range("M64").Value = -2128.5
$ws.Range("N64").Value = -15340.0005
$ws.Range("H67").Value = 4181.727
$ws.Range("I67").Value = 799.5
$ws.Range("J67").Value = 4933.3335
$ws.Range("K67").Value = 2398.5
$ws.Range("L67").Value = 14800.0005
$ws.Range("M67").Value = -1462.5
$ws.Range("N67").Value = -16672.0005
$ws.Range("H82").Value = 10500.857
$ws.Range("I82").Value = 2006.5
$ws.Range("J82").Value = 11916.583
$ws.Range("K82").Value = 6019.5
$ws.Range("L82").Value = 35749.749
$ws.Range("M82").Value = -5613.5
$ws.Range("N82").Value = -36561.749
$ws.Range("H85").Value = 10500.857
$ws.Range("I85").Value = 2006.5
$ws.Range("J85").Value = 11916.583
$ws.Range("K85").Value = 6019.5
$ws.Range("L85").Value = 35749.749
$ws.Range("M85").Value = -4615.5
$ws.Range("N85").Value = -38557.749
$ws.Range("H95").Value = 6587.8335
$ws.Range("J95").Value = 6587.8335
$ws.Range("L95").Value = 19763.5005
$ws.Range("N95").Value = -23881.5005
$ws.Range("H122").Value = 759.8182
$ws.Range("I122").Value = 701.2857
$ws.Range("J122").Value = 862.25
$ws.Range("K122").Value = 6311.571300000001
$ws.Range("L122").Value = 7760.25
$ws.Range("M122").Value = -3861.571300000001
$ws.Range("N122").Value = -12660.25
$ws.Range("H131").Value = 19232112
$ws.Range("I131").Value = 100000370
$ws.Range("J131").Value = 1574.4286
$ws.Range("K131").Value = 300001110
$ws.Range("L131").Value = 4723.2858
$ws.Range("M131").Value = -299996070
$ws.Range("N131").Value = -14803.2858
$ws.Range("H134").Value = 4941.485
$ws.Range("I134").Value = 2434.75
$ws.Range("J134").Value = 5743.64
$ws.Range("K134").Value = 7304.25
$ws.Range("L134").Value = 17230.92
$ws.Range("M134").Value = -2234.25
$ws.Range("N134").Value = -27370.92
$ws.Range("H135").Value = 1932.0625
$ws.Range("I135").Value = 2064.8572
$ws.Range("K135").Value = 18583.7148
$ws.Range("M135").Value = -16048.7148

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1344.2858
$ws.Range("I22").Value = 1194.5454
$ws.Range("J22").Value = 1893.3334
$ws.Range("K22").Value = 1194.5454
$ws.Range("L22").Value = 1893.3334
$ws.Range("M22").Value = -899.5454
$ws.Range("N22").Value = -2483.3334
$ws.Range("H27").Value = 1344.2858
$ws.Range("I27").Value = 1194.5454
$ws.Range("J27").Value = 1893.3334
$ws.Range("K27").Value = 1194.5454
$ws.Range("L27").Value = 1893.3334
$ws.Range("M27").Value = -1087.5454
$ws.Range("N27").Value = -2107.3334
$ws.Range("H46").Value = 7950
$ws.Range("I46").Value = 1300.5
$ws.Range("J46").Value = 9279.9
$ws.Range("K46").Value = 1300.5
$ws.Range("L46").Value = 9279.9
$ws.Range("M46").Value = -1112.5
$ws.Range("N46").Value = -9655.9
$ws.Range("H68").Value = 1301.5883
$ws.Range("I68").Value = 1131.75
$ws.Range("J68").Value = 1709.2
$ws.Range("K68").Value = 1131.75
$ws.Range("L68").Value = 1709.2
$ws.Range("M68").Value = -382.75
$ws.Range("N68").Value = -3207.2
$ws.Range("H71").Value = 1301.5883
$ws.Range("I71").Value = 1131.75
$ws.Range("J71").Value = 1709.2
$ws.Range("K71").Value = 5658.75
$ws.Range("L71").Value = 8546
$ws.Range("M71").Value = -1914.75
$ws.Range("N71").Value = -16034

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2120
$ws.Range("J4").Value = 2120
$ws.Range("L4").Value = 2120
$ws.Range("N4").Value = -2346
$ws.Range("H122").Value = 8969584
$ws.Range("I122").Value = 9633812
$ws.Range("J122").Value = 2502.5
$ws.Range("K122").Value = 28901436
$ws.Range("L122").Value = 7507.5
$ws.Range("M122").Value = -28898986
$ws.Range("N122").Value = -12407.5
$ws.Range("H136").Value = 630.9091
$ws.Range("I136").Value = 304.44446
$ws.Range("K136").Value = 913.33338
$ws.Range("M136").Value = 1636.66662
